$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-24 Wednesday", "2024-07-25 Thursday"),
    @("139÷6=", "675÷9="),
    @("344÷6=", "422÷3="),
    @("252÷2=", "815÷3="),
    @("393÷2=", "482÷7="),
    @("392÷2=", "132÷8="),
    @("894÷7=", "719÷2="),
    @("521÷4=", "520÷6="),
    @("689÷3=", "399÷2="),
    @("352÷2=", "114÷3="),
    @("159÷5=", "570÷9="),
    @("523÷9=", "861÷9="),
    @("420÷9=", "595÷3="),
    @("519÷3=", "519÷5="),
    @("324÷8=", "219÷4="),
    @("417÷6=", "609÷9="),
    @("136÷2=", "866÷4="),
    @("288÷6=", "350÷4="),
    @("138÷8=", "823÷7="),
    @("658÷2=", "149÷3="),
    @("665÷5=", "194÷9="),
    @("505÷4=", "719÷3="),
    @("114÷2=", "922÷8="),
    @("874÷6=", "963÷7="),
    @("174÷3=", "388÷9="),
    @("124÷9=", "231÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
